$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new "Random" sample row (row 3) with the same shape as row 2
$ws.Range("A3").Value = 42600.881238425929
$ws.Range("B3").Value = "Random"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 42
$ws.Range("I3").Value = 58
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 60
$ws.Range("M3").Value = 40

# Widen column A slightly to fit the new content (bestFit column A: 13.85546875 -> 14.85546875)
$ws.Columns.Item(1).ColumnWidth = 14
